$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("n1_d40")
$ws.Range("F2").Value = 0.118142
$ws.Range("G2").Value = 10.5
$ws.Range("F3").Value = 0.118242
$ws.Range("G3").Value = 10.8
$ws.Range("F4").Value = 0.119069
$ws.Range("G4").Value = 10.5
$ws.Range("F5").Value = 0.119191
$ws.Range("G5").Value = 10.5
$ws.Range("F6").Value = 0.119295
$ws.Range("G6").Value = 10.5
$ws.Range("F7").Value = 0.118984
$ws.Range("G7").Value = 10.4
$ws.Range("F8").Value = 0.120257
$ws.Range("G8").Value = 10.8
$ws.Range("F9").Value = 0.119275
$ws.Range("G9").Value = 10.5
$ws.Range("F10").Value = 0.119251
$ws.Range("G10").Value = 10.8
$ws.Range("F11").Value = 0.120586
$ws.Range("G11").Value = 10.7
$ws.Range("F12").Value = 0.1192292
$ws.Range("G12").Value = 10.6

$ws = $wb.Worksheets.Item("n1_d60")
$ws.Range("F2").Value = 0.156424
$ws.Range("G2").Value = 16.7
$ws.Range("F3").Value = 0.155539
$ws.Range("G3").Value = 16.7
$ws.Range("F4").Value = 0.155525
$ws.Range("G4").Value = 16.7
$ws.Range("F5").Value = 0.15739
$ws.Range("G5").Value = 16.7
$ws.Range("F6").Value = 0.155582
$ws.Range("G6").Value = 16.7
$ws.Range("F7").Value = 0.16161
$ws.Range("G7").Value = 16.7
$ws.Range("F8").Value = 0.157626
$ws.Range("G8").Value = 16.7
$ws.Range("F9").Value = 0.155436
$ws.Range("G9").Value = 16.7
$ws.Range("F10").Value = 0.155499
$ws.Range("G10").Value = 16.7
$ws.Range("F11").Value = 0.156513
$ws.Range("G11").Value = 16.7
$ws.Range("F12").Value = 0.1567144
$ws.Range("G12").Value = 16.7

$ws = $wb.Worksheets.Item("n1_d80")
$ws.Range("F2").Value = 0.191801
$ws.Range("G2").Value = 24.2
$ws.Range("F3").Value = 0.193131
$ws.Range("G3").Value = 24.2
$ws.Range("F4").Value = 0.19076
$ws.Range("G4").Value = 24.8
$ws.Range("F5").Value = 0.190606
$ws.Range("G5").Value = 24.7
$ws.Range("F6").Value = 0.190635
$ws.Range("G6").Value = 24.2
$ws.Range("F7").Value = 0.19826
$ws.Range("G7").Value = 24.2
$ws.Range("F8").Value = 0.195458
$ws.Range("G8").Value = 24.1
$ws.Range("F9").Value = 0.193167
$ws.Range("G9").Value = 24.2
$ws.Range("F10").Value = 0.194254
$ws.Range("G10").Value = 24.2
$ws.Range("F11").Value = 0.192978
$ws.Range("G11").Value = 24.2
$ws.Range("F12").Value = 0.193105
$ws.Range("G12").Value = 24.3

$ws = $wb.Worksheets.Item("n1_d100")
$ws.Range("F2").Value = 0.222794
$ws.Range("G2").Value = 29.7
$ws.Range("F3").Value = 0.223125
$ws.Range("G3").Value = 28.7
$ws.Range("F4").Value = 0.224586
$ws.Range("G4").Value = 29.1
$ws.Range("F5").Value = 0.224769
$ws.Range("G5").Value = 29.9
$ws.Range("F6").Value = 0.22386
$ws.Range("G6").Value = 29.1
$ws.Range("F7").Value = 0.227626
$ws.Range("G7").Value = 29.3
$ws.Range("F8").Value = 0.226987
$ws.Range("G8").Value = 29.2
$ws.Range("F9").Value = 0.223776
$ws.Range("G9").Value = 29.3
$ws.Range("F10").Value = 0.229073
$ws.Range("G10").Value = 29.3
$ws.Range("F11").Value = 0.230564
$ws.Range("G11").Value = 28.7
$ws.Range("F12").Value = 0.225716
$ws.Range("G12").Value = 29.23
